$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) Paragraph 3: "Predicting arms transfers with using multiple
#    data sets; " -> Turkish replacement text, split per-word with
#    proofErr spell-check wrappers (to match Word's auto-flagging
#    of non-dictionary words).
# ---------------------------------------------------------------
$p3 = $d.Paragraphs(3).Range
$frag3 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="027450E1" w14:textId="45531412" w:rsidR="00AB6810" w:rsidRDefault="00AB6810" w:rsidP="0086015E"><w:pPr><w:jc w:val="both"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Buradaki</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dosyayi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>degistirmemiz</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>gerekiyor</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p3.InsertXML($frag3)

# ---------------------------------------------------------------
# 2) Paragraph 6 ("It is a fact..."): insert a collapsed _GoBack
#    bookmark between "ha" and "s" of "has become a globalized
#    problem" (last-edit-position marker Word stamps on save).
# ---------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("warfare ha", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pos = $rng.End
$bmRange = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------
# 3) Paragraph 7 ("Stockholm International..."): merge 3 runs into 1.
# ---------------------------------------------------------------
$p7 = $d.Paragraphs(7).Range
$frag7 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="499F310A" w14:textId="39209C28" w:rsidR="000D32B0" w:rsidRDefault="006D3067" w:rsidP="006D3067"><w:pPr><w:spacing w:after="120"/><w:ind w:firstLine="288"/><w:jc w:val="both"/></w:pPr><w:r><w:t>Stockholm International Peace Research Institute (SIPRI) database is a unique resource for researchers, policy-makers and analysts, the media and civil society interested in monitoring and measuring the international flow of major conventional arms.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p7.InsertXML($frag7)

# ---------------------------------------------------------------
# 4) Paragraph 11 ("Main purpose..."): merge the middle two runs
#    (" with prediction of " + "flow of major conventional arms").
# ---------------------------------------------------------------
$p11 = $d.Paragraphs(11).Range
$frag11 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="79CC9B76" w14:textId="17BB58EE" w:rsidR="000D32B0" w:rsidRDefault="000D32B0" w:rsidP="00552EF5"><w:pPr><w:spacing w:after="120"/><w:ind w:firstLine="288"/><w:jc w:val="both"/></w:pPr><w:r><w:t>Main purpose is to make a contribution to estimation phase of a security crisis</w:t></w:r><w:r w:rsidR="006D3067"><w:t xml:space="preserve"> with prediction of flow of major conventional arms</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p11.InsertXML($frag11)

# ---------------------------------------------------------------
# 5) Paragraph 17 ("Can we predict the future arms transfers?"):
#    merge "future "+"arms " and "t"+"ransfers" runs.
# ---------------------------------------------------------------
$p17 = $d.Paragraphs(17).Range
$frag17 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="5A58ABB6" w14:textId="7DED7BA3" w:rsidR="0005134F" w:rsidRPr="00DB4834" w:rsidRDefault="0005134F" w:rsidP="0005134F"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="120"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr></w:pPr><w:r w:rsidRPr="00DB4834"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">Can we predict the </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">future arms </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>transfers</w:t></w:r><w:r w:rsidRPr="00DB4834"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>?</w:t></w:r><w:r w:rsidR="00E0702D"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>‘</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidR="00E0702D"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:vertAlign w:val="subscript"/></w:rPr><w:sym w:font="Symbol" w:char="F0B9"/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p17.InsertXML($frag17)

# ---------------------------------------------------------------
# 6) Paragraph 22 ("Importer/exporter the total trend-indicator
#    value (TIV) tables"): merge 5 runs into 1.
# ---------------------------------------------------------------
$p22 = $d.Paragraphs(22).Range
$frag22 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="6324A445" w14:textId="4F0C5FA1" w:rsidR="00472ADD" w:rsidRDefault="00472ADD" w:rsidP="00472ADD"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:spacing w:after="120"/><w:jc w:val="both"/></w:pPr><w:r><w:t>Importer/exporter the total trend-indicator value (TIV) tables</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p22.InsertXML($frag22)

Write-Output "edits applied"
